# The underlying diff shows the "Rules" worksheet's cell C10 (row for "R30")
# changing its value from 18 to 1. All other differences in the XML diff
# (numFmts, xfId/collapsed attributes, col attribute ordering) are
# serialization artifacts of the external tool that produced the diff and
# are not semantic workbook changes, so only the cell value is updated here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C10").Value = 1
